$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 78.333336
$ws.Cells.Item(6, 9).Value = 91.5
$ws.Cells.Item(6, 10).Value = 12.5
$ws.Cells.Item(6, 11).Value = 274.5
$ws.Cells.Item(6, 12).Value = 37.5
$ws.Cells.Item(6, 13).Value = -162.5
$ws.Cells.Item(6, 14).Value = -261.5
$ws.Cells.Item(48, 8).Value = 2998.6667
$ws.Cells.Item(48, 10).Value = 2998.6667
$ws.Cells.Item(48, 12).Value = 8996.000100000001
$ws.Cells.Item(48, 14).Value = -9580.000100000001
$ws.Cells.Item(56, 8).Value = 2998.6667
$ws.Cells.Item(56, 10).Value = 2998.6667
$ws.Cells.Item(56, 12).Value = 8996.000100000001
$ws.Cells.Item(56, 14).Value = -10064.0001
$ws.Cells.Item(110, 8).Value = 40200.8
$ws.Cells.Item(110, 10).Value = 40200.8
$ws.Cells.Item(110, 12).Value = 40200.8
$ws.Cells.Item(110, 14).Value = -48380.8
$ws.Cells.Item(137, 8).Value = 12475.968
$ws.Cells.Item(137, 9).Value = 2362.65
$ws.Cells.Item(137, 10).Value = 30863.818
$ws.Cells.Item(137, 11).Value = 7087.950000000001
$ws.Cells.Item(137, 12).Value = 92591.454
$ws.Cells.Item(137, 13).Value = -4537.950000000001
$ws.Cells.Item(137, 14).Value = -97691.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 6373.636
$ws.Cells.Item(2, 9).Value = 1851.8334
$ws.Cells.Item(2, 10).Value = 11799.8
$ws.Cells.Item(2, 11).Value = 1851.8334
$ws.Cells.Item(2, 12).Value = 11799.8
$ws.Cells.Item(2, 13).Value = -1738.8334
$ws.Cells.Item(2, 14).Value = -12025.8
$ws.Cells.Item(32, 8).Value = 10577.444
$ws.Cells.Item(32, 9).Value = 3888.9524
$ws.Cells.Item(32, 10).Value = 26183.926
$ws.Cells.Item(32, 11).Value = 3888.9524
$ws.Cells.Item(32, 12).Value = 26183.926
$ws.Cells.Item(32, 13).Value = -3601.9524
$ws.Cells.Item(32, 14).Value = -26757.926
$ws.Cells.Item(63, 8).Value = 5805
$ws.Cells.Item(63, 9).Value = 6189.1665
$ws.Cells.Item(63, 10).Value = 3500
$ws.Cells.Item(63, 11).Value = 6189.1665
$ws.Cells.Item(63, 12).Value = 3500
$ws.Cells.Item(63, 13).Value = -5503.1665
$ws.Cells.Item(63, 14).Value = -4872
$ws.Cells.Item(66, 8).Value = 5805
$ws.Cells.Item(66, 9).Value = 6189.1665
$ws.Cells.Item(66, 10).Value = 3500
$ws.Cells.Item(66, 11).Value = 30945.8325
$ws.Cells.Item(66, 12).Value = 17500
$ws.Cells.Item(66, 13).Value = -27513.8325
$ws.Cells.Item(66, 14).Value = -24364
$ws.Cells.Item(102, 8).Value = 10968.5
$ws.Cells.Item(102, 9).Value = 2783.8635
$ws.Cells.Item(102, 11).Value = 2783.8635
$ws.Cells.Item(102, 13).Value = -1161.8635
$ws.Cells.Item(116, 8).Value = 6373.636
$ws.Cells.Item(116, 9).Value = 1851.8334
$ws.Cells.Item(116, 10).Value = 11799.8
$ws.Cells.Item(116, 11).Value = 1851.8334
$ws.Cells.Item(116, 12).Value = 11799.8
$ws.Cells.Item(116, 13).Value = 442.1666
$ws.Cells.Item(116, 14).Value = -16387.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 6373.636
$ws.Cells.Item(3, 9).Value = 1851.8334
$ws.Cells.Item(3, 10).Value = 11799.8
$ws.Cells.Item(3, 11).Value = 1851.8334
$ws.Cells.Item(3, 12).Value = 11799.8
$ws.Cells.Item(3, 13).Value = -1737.8334
$ws.Cells.Item(3, 14).Value = -12027.8
$ws.Cells.Item(105, 8).Value = 1978.375
$ws.Cells.Item(105, 9).Value = 1024.5714
$ws.Cells.Item(105, 11).Value = 1024.5714
$ws.Cells.Item(105, 13).Value = 722.4286
$ws.Cells.Item(107, 8).Value = 2535.348
$ws.Cells.Item(107, 9).Value = 2204.3125
$ws.Cells.Item(107, 11).Value = 2204.3125
$ws.Cells.Item(107, 13).Value = -284.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 151.57692
$ws.Cells.Item(7, 9).Value = 55.35294
$ws.Cells.Item(7, 10).Value = 333.33334
$ws.Cells.Item(7, 11).Value = 55.35294
$ws.Cells.Item(7, 12).Value = 333.33334
$ws.Cells.Item(7, 13).Value = 57.64706
$ws.Cells.Item(7, 14).Value = -559.33334
$ws.Cells.Item(16, 8).Value = 5906.3076
$ws.Cells.Item(16, 9).Value = 2097.875
$ws.Cells.Item(16, 10).Value = 11999.8
$ws.Cells.Item(16, 11).Value = 2097.875
$ws.Cells.Item(16, 12).Value = 11999.8
$ws.Cells.Item(16, 13).Value = -1810.875
$ws.Cells.Item(16, 14).Value = -12573.8
$ws.Cells.Item(41, 8).Value = 7621.8335
$ws.Cells.Item(41, 9).Value = 7621.8335
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 7621.8335
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = -7193.8335
$ws.Cells.Item(41, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 5906.3076
$ws.Cells.Item(113, 9).Value = 2097.875
$ws.Cells.Item(113, 10).Value = 11999.8
$ws.Cells.Item(113, 11).Value = 2097.875
$ws.Cells.Item(113, 12).Value = 11999.8
$ws.Cells.Item(113, 13).Value = 72.125
$ws.Cells.Item(113, 14).Value = -16339.8
$ws.Cells.Item(134, 8).Value = 26321476
$ws.Cells.Item(134, 9).Value = 1059.3478
$ws.Cells.Item(134, 10).Value = 66679450
$ws.Cells.Item(134, 11).Value = 3178.0434
$ws.Cells.Item(134, 12).Value = 200038350
$ws.Cells.Item(134, 13).Value = -643.0434
$ws.Cells.Item(134, 14).Value = -200043420

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 37
$ws.Cells.Item(12, 10).Value = 1.5
$ws.Cells.Item(12, 12).Value = 4.5
$ws.Cells.Item(12, 14).Value = -350.5
$ws.Cells.Item(26, 8).Value = 1630.25
$ws.Cells.Item(26, 9).Value = 2587
$ws.Cells.Item(26, 10).Value = 673.5
$ws.Cells.Item(26, 11).Value = 7761
$ws.Cells.Item(26, 12).Value = 2020.5
$ws.Cells.Item(26, 13).Value = -7473
$ws.Cells.Item(26, 14).Value = -2596.5
$ws.Cells.Item(97, 8).Value = 7920
$ws.Cells.Item(97, 9).Value = 466.66666
$ws.Cells.Item(97, 10).Value = 19100
$ws.Cells.Item(97, 11).Value = 1399.99998
$ws.Cells.Item(97, 12).Value = 57300
$ws.Cells.Item(97, 13).Value = -903.9999800000001
$ws.Cells.Item(97, 14).Value = -58292
$ws.Cells.Item(116, 8).Value = 3619
$ws.Cells.Item(116, 9).Value = 3943.25
$ws.Cells.Item(116, 10).Value = 1025
$ws.Cells.Item(116, 11).Value = 11829.75
$ws.Cells.Item(116, 12).Value = 3075
$ws.Cells.Item(116, 13).Value = -8387.75
$ws.Cells.Item(116, 14).Value = -9959
$ws.Cells.Item(136, 8).Value = 2000
$ws.Cells.Item(136, 9).Value = 2000
$ws.Cells.Item(136, 11).Value = 6000
$ws.Cells.Item(136, 13).Value = -900

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 117.833336
$ws.Cells.Item(12, 9).Value = 100
$ws.Cells.Item(12, 11).Value = 100
$ws.Cells.Item(12, 13).Value = 40
$ws.Cells.Item(18, 8).Value = 1505333.1
$ws.Cells.Item(18, 10).Value = 5999.75
$ws.Cells.Item(18, 12).Value = 5999.75
$ws.Cells.Item(18, 14).Value = -6585.75
$ws.Cells.Item(21, 8).Value = 1609600
$ws.Cells.Item(21, 9).Value = 4005000
$ws.Cells.Item(21, 10).Value = 12666.667
$ws.Cells.Item(21, 11).Value = 4005000
$ws.Cells.Item(21, 12).Value = 12666.667
$ws.Cells.Item(21, 13).Value = -4004827
$ws.Cells.Item(21, 14).Value = -13012.667
$ws.Cells.Item(29, 8).Value = 5797
$ws.Cells.Item(29, 9).Value = 5797
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 5797
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = -5507
$ws.Cells.Item(29, 14).ClearContents()
$ws.Cells.Item(30, 8).Value = 1609600
$ws.Cells.Item(30, 9).Value = 4005000
$ws.Cells.Item(30, 10).Value = 12666.667
$ws.Cells.Item(30, 11).Value = 4005000
$ws.Cells.Item(30, 12).Value = 12666.667
$ws.Cells.Item(30, 13).Value = -4004895
$ws.Cells.Item(30, 14).Value = -12876.667
$ws.Cells.Item(64, 8).Value = 49900
$ws.Cells.Item(64, 10).Value = 49900
$ws.Cells.Item(64, 12).Value = 49900
$ws.Cells.Item(64, 14).Value = -50396
$ws.Cells.Item(67, 8).Value = 49900
$ws.Cells.Item(67, 10).Value = 49900
$ws.Cells.Item(67, 12).Value = 49900
$ws.Cells.Item(67, 14).Value = -51616
$ws.Cells.Item(132, 8).Value = 9254.77
$ws.Cells.Item(132, 9).Value = 9346.695
$ws.Cells.Item(132, 10).Value = 8550
$ws.Cells.Item(132, 11).Value = 28040.085
$ws.Cells.Item(132, 12).Value = 25650
$ws.Cells.Item(132, 13).Value = -25510.085
$ws.Cells.Item(132, 14).Value = -30710

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(30, 8).Value = 1859.1666
$ws.Cells.Item(30, 9).Value = 288.75
$ws.Cells.Item(30, 10).Value = 5000
$ws.Cells.Item(30, 11).Value = 288.75
$ws.Cells.Item(30, 12).Value = 5000
$ws.Cells.Item(30, 13).Value = -180.75
$ws.Cells.Item(30, 14).Value = -5216
$ws.Cells.Item(31, 8).Value = 1250
$ws.Cells.Item(31, 9).Value = 1250
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 1250
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -1002
$ws.Cells.Item(31, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 9161.565000000001
$ws.Cells.Item(122, 9).Value = 7027.25
$ws.Cells.Item(122, 11).Value = 21081.75
$ws.Cells.Item(122, 13).Value = -18631.75
$ws.Cells.Item(132, 8).Value = 1490849.5
$ws.Cells.Item(132, 9).Value = 5203.385
$ws.Cells.Item(132, 11).Value = 15610.155
$ws.Cells.Item(132, 13).Value = -13080.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 2347049.2
$ws.Cells.Item(4, 9).Value = 10000325
$ws.Cells.Item(4, 10).Value = 51066.4
$ws.Cells.Item(4, 11).Value = 10000325
$ws.Cells.Item(4, 12).Value = 51066.4
$ws.Cells.Item(4, 13).Value = -10000212
$ws.Cells.Item(4, 14).Value = -51292.4
$ws.Cells.Item(69, 8).Value = 44577
$ws.Cells.Item(69, 10).Value = 47334.4
$ws.Cells.Item(69, 12).Value = 47334.4
$ws.Cells.Item(69, 14).Value = -48832.4
$ws.Cells.Item(72, 8).Value = 44577
$ws.Cells.Item(72, 10).Value = 47334.4
$ws.Cells.Item(72, 12).Value = 142003.2
$ws.Cells.Item(72, 14).Value = -149491.2
$ws.Cells.Item(107, 8).Value = 6299.4
$ws.Cells.Item(107, 9).Value = 750
$ws.Cells.Item(107, 11).Value = 2250
$ws.Cells.Item(107, 13).Value = -330
$ws.Cells.Item(126, 8).Value = 6785
$ws.Cells.Item(126, 10).Value = 25748.5
$ws.Cells.Item(126, 12).Value = 77245.5
$ws.Cells.Item(126, 14).Value = -82185.5
$ws.Cells.Item(132, 8).Value = 7543.8057
$ws.Cells.Item(132, 9).Value = 1739.9546
$ws.Cells.Item(132, 10).Value = 16664.143
$ws.Cells.Item(132, 11).Value = 5219.8638
$ws.Cells.Item(132, 12).Value = 49992.429
$ws.Cells.Item(132, 13).Value = -2689.8638
$ws.Cells.Item(132, 14).Value = -55052.429
